# Insert a new "債務" (debt) worksheet right before the "事業投資" sheet,
# matching the target xlsx produced by the commit "fix some dirty files".

$wb = $excel.ActiveWorkbook

# Add the new sheet immediately before "事業投資"
$investSheet = $wb.Worksheets.Item("事業投資")
$debtSheet = $wb.Worksheets.Add($investSheet, $null)
$debtSheet.Name = "債務"

# Mirror the look of the other property sheets: bold+bordered header row
# (style used on row 1 / column A of every other sheet in this workbook),
# reused here via copy/paste-format from the "保險" sheet instead of being
# rebuilt by hand.
$fmtSrc = $wb.Worksheets.Item("保險")
$fmtSrc.Range("B1:K1").Copy()
$debtSheet.Range("B1:N1").PasteSpecial(-4122)   # xlPasteFormats
$fmtSrc.Range("A2").Copy()
$debtSheet.Range("A2:A3").PasteSpecial(-4122)
$fmtSrc.Range("B2:K2").Copy()
$debtSheet.Range("B2:N3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "date" column (J) data cells hold a value that looks like a date
# ("2012-04-18"); force them to Text first so the value is kept verbatim
# instead of being auto-converted to a date serial number. (J1 is just the
# header label "date" and does not need this.)
$debtSheet.Range("J2:J3").NumberFormat = "@"

# Header row (row 1) - column A is left blank, same as the other sheets
$debtSheet.Range("B1").Value = "species"
$debtSheet.Range("C1").Value = "debtor"
$debtSheet.Range("D1").Value = "owner"
$debtSheet.Range("E1").Value = "total"
$debtSheet.Range("F1").Value = "register_date"
$debtSheet.Range("G1").Value = "register_reason"
$debtSheet.Range("H1").Value = "property_category"
$debtSheet.Range("I1").Value = "category"
$debtSheet.Range("J1").Value = "date"
$debtSheet.Range("K1").Value = "legislator_name"
$debtSheet.Range("L1").Value = "legislator_id"
$debtSheet.Range("M1").Value = "source_file"
$debtSheet.Range("N1").Value = "index"

# Row 2
$debtSheet.Range("A2").Value = 111
$debtSheet.Range("B2").Value = "房屋貸款"
$debtSheet.Range("C2").Value = "邱議瑩"
$debtSheet.Range("D2").Value = "京城銀行忠孝分行臺北市南港區忠孝東路"
$debtSheet.Range("E2").Value = 8245855
$debtSheet.Range("F2").Value = "98年11月16日"
$debtSheet.Range("G2").Value = "購屋"
$debtSheet.Range("H2").Value = "debt"
$debtSheet.Range("I2").Value = "normal"
$debtSheet.Range("J2").Value = "2012-04-18"
$debtSheet.Range("K2").Value = "邱議瑩"
$debtSheet.Range("L2").Value = 913
$debtSheet.Range("M2").Value = "tmped121"
$debtSheet.Range("N2").Value = 111

# Row 3
$debtSheet.Range("A3").Value = 112
$debtSheet.Range("B3").Value = "房屋貸款"
$debtSheet.Range("C3").Value = "邱議瑩"
$debtSheet.Range("D3").Value = "台灣銀行高雄三民分行高雄市三民區九如二路"
$debtSheet.Range("E3").Value = 10000000
$debtSheet.Range("F3").Value = "99年06月21曰"
$debtSheet.Range("G3").Value = "購屋"
$debtSheet.Range("H3").Value = "debt"
$debtSheet.Range("I3").Value = "normal"
$debtSheet.Range("J3").Value = "2012-04-18"
$debtSheet.Range("K3").Value = "邱議瑩"
$debtSheet.Range("L3").Value = 913
$debtSheet.Range("M3").Value = "tmped121"
$debtSheet.Range("N3").Value = 112
